$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 11796
$ws1.Range("F3").Value = 11498
$ws1.Range("F4").Value = 609
$ws1.Range("F6").Value = 1038
$ws1.Range("F11").Value = 10833
$ws1.Range("F12").Value = 4185
$ws1.Range("F13").Value = 19
$ws1.Range("F14").Value = 13
$ws1.Range("F17").Value = 1055
$ws1.Range("F18").Value = 56
$ws1.Range("F19").Value = 3
$ws1.Range("F22").Value = 11155
$ws1.Range("F23").Value = 10946
$ws1.Range("F28").Value = 30

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 11796
$ws4.Range("F3").Value = 11498
$ws4.Range("F4").Value = 609
$ws4.Range("F6").Value = 1038
$ws4.Range("F11").Value = 10833
$ws4.Range("F12").Value = 4185
$ws4.Range("F13").Value = 19
$ws4.Range("F14").Value = 13
$ws4.Range("F17").Value = 1055
$ws4.Range("F18").Value = 56
$ws4.Range("F19").Value = 3
$ws4.Range("F22").Value = 11155
$ws4.Range("F23").Value = 10946
$ws4.Range("F28").Value = 30
